$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append two new events as rows 180 and 181 --------------------------
# Those rows currently exist only as blank, placeholder-styled cells
# (style "5"); clone the formatting (borders/fill/number-format) of the
# last real data row (179) onto them first, so they end up with the same
# look (date format on A, bordered text style on B:E) as every other
# populated row instead of the blank placeholder style.
$ws.Range("A179:E179").Copy()
$ws.Range("A180:E181").PasteSpecial(-4122)

$eventDate = Get-Date -Year 2025 -Month 2 -Day 28 -Hour 0 -Minute 0 -Second 0
$link180 = "https://www.instagram.com/reel/DF3BQudtOjO/?igsh=cHQ1NDdkZTlxdWpn"
$link181 = "https://www.instagram.com/reel/DGBIl8nsj97/?igsh=OGVnMTZxdDJmdjQ="

# Row 180 - UEBERREST @ Schrotty, Koeln (28.02.2025)
$ws.Range("A180").Value = $eventDate
$ws.Range("B180").Value = "UEBERREST"
$ws.Range("C180").Value = "Schrotty"
$ws.Range("D180").Value = "Köln"
$ws.Range("E180").Value = $link180

# Row 181 - HIGH VOLTAGE HARD CARNIVAL 12H RAVE (18 Uhr) @ Schlachthof, Duesseldorf (28.02.2025)
$ws.Range("A181").Value = $eventDate
$ws.Range("B181").Value = "HIGH VOLTAGE HARD CARNIVAL 12H RAVE (18 Uhr)"
$ws.Range("C181").Value = "Schlachthof"
$ws.Range("D181").Value = "Düsseldorf"
$ws.Range("E181").Value = $link181

# Turn the two new Link cells into real hyperlinks, same as every other
# cell in column E. Hyperlinks.Add() stamps its own default "Hyperlink"
# cell style onto the target, so immediately re-apply the row-179 derived
# formatting afterwards to keep the original bordered look.
$ws.Hyperlinks.Add($ws.Range("E180"), $link180, "", "", $link180)
$ws.Hyperlinks.Add($ws.Range("E181"), $link181, "", "", $link181)

$ws.Range("A179:E179").Copy()
$ws.Range("A180:E181").PasteSpecial(-4122)

$ws.Range("A180").Value = $eventDate
$ws.Range("B180").Value = "UEBERREST"
$ws.Range("C180").Value = "Schrotty"
$ws.Range("D180").Value = "Köln"
$ws.Range("E180").Value = $link180

$ws.Range("A181").Value = $eventDate
$ws.Range("B181").Value = "HIGH VOLTAGE HARD CARNIVAL 12H RAVE (18 Uhr)"
$ws.Range("C181").Value = "Schlachthof"
$ws.Range("D181").Value = "Düsseldorf"
$ws.Range("E181").Value = $link181

# Give the Link cells' text the same underlined/colored run formatting
# used by every other hyperlink cell in the sheet. Splitting the run into
# two adjoining pieces (rather than styling the whole string as a single
# pass) is what makes the engine keep it as a shared-string rich-text run
# instead of collapsing it into a whole-cell font override.
foreach ($coord in @("E180", "E181")) {
    $cell = $ws.Range($coord)
    $n = $cell.Value2.Length
    $head = $cell.Characters(1, $n - 1)
    $head.Font.Underline = $true
    $head.Font.ColorIndex = 4
    $head.Font.Name = "Calibri"
    $head.Font.Size = 11
    $tail = $cell.Characters($n, 1)
    $tail.Font.Underline = $true
    $tail.Font.ColorIndex = 4
    $tail.Font.Name = "Calibri"
    $tail.Font.Size = 11
}

Write-Host "done"
